# Update "想去人数" (interested-people count) figures on the
# "展览" (Exhibition) and "全部类型" (All types) sheets.
#
# 展览 sheet:
#   F2: 1374 -> 1377
#   F3: 2911 -> 2917
#   F4: 13   -> 15
#
# 全部类型 sheet:
#   F3: 1374 -> 1377
#   F4: 2911 -> 2917
#   F5: 13   -> 15

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1377
$ws1.Range("F3").Value = 2917
$ws1.Range("F4").Value = 15

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1377
$ws4.Range("F4").Value = 2917
$ws4.Range("F5").Value = 15
